$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.56542157841475
$ws.Range("D2").Value = 8.990910777557627
$ws.Range("E2").Value = 14.90835905290406
$ws.Range("F2").Value = 51.60254544065107
$ws.Range("G2").Value = 3.594859047301906
$ws.Range("I2").Value = 19.56988205597989
$ws.Range("M2").Value = 31.99538715734932
$ws.Range("B3").Value = 8.231670994341819
$ws.Range("D3").Value = 8.712696350048354
$ws.Range("E3").Value = 14.08648243838648
$ws.Range("F3").Value = 49.50204601729268
$ws.Range("G3").Value = 3.607827519654669
$ws.Range("I3").Value = 19.41091534122405
$ws.Range("M3").Value = 30.44997188032105
$ws.Range("B4").Value = 8.019168844761055
$ws.Range("D4").Value = 8.545499728442572
$ws.Range("E4").Value = 13.55594164246932
$ws.Range("F4").Value = 48.20784397318952
$ws.Range("G4").Value = 3.616113791610688
$ws.Range("I4").Value = 19.32002503466039
$ws.Range("M4").Value = 29.46983455025978
$ws.Range("B5").Value = 7.930764315880739
$ws.Range("D5").Value = 8.478378044973393
$ws.Range("E5").Value = 13.33334095594995
$ws.Range("F5").Value = 47.68020151474055
$ws.Range("G5").Value = 3.619572984196569
$ws.Range("I5").Value = 19.2847232762369
$ws.Range("M5").Value = 29.06305886772537
$ws.Range("B6").Value = 7.915978690211833
$ws.Range("D6").Value = 8.46729660369035
$ws.Range("E6").Value = 13.29599551921294
$ws.Range("F6").Value = 47.59259939249599
$ws.Range("G6").Value = 3.620152393142308
$ws.Range("I6").Value = 19.27896769174965
$ws.Range("M6").Value = 28.9950854851875
$ws.Range("B7").Value = 8.017983773905787
$ws.Range("D7").Value = 8.544590270910465
$ws.Range("E7").Value = 13.55296531998253
$ws.Range("F7").Value = 48.20072776995811
$ws.Range("G7").Value = 3.616160108135272
$ws.Range("I7").Value = 19.3195418492189
$ws.Range("M7").Value = 29.46437772572614
$ws.Range("B8").Value = 8.451971825018497
$ws.Range("D8").Value = 8.894289112148462
$ws.Range("E8").Value = 14.63041239056321
$ws.Range("F8").Value = 50.87973058789609
$ws.Range("G8").Value = 3.599264092986509
$ws.Range("I8").Value = 19.51370608116482
$ws.Range("M8").Value = 31.46928617786296
$ws.Range("B9").Value = 9.239219066721178
$ws.Range("D9").Value = 9.60454587072288
$ws.Range("E9").Value = 16.53402456961178
$ws.Range("F9").Value = 56.06298538161366
$ws.Range("G9").Value = 3.568646813689276
$ws.Range("I9").Value = 19.94551427612702
$ws.Range("M9").Value = 35.13506105074688
$ws.Range("B10").Value = 9.774561620633685
$ws.Range("D10").Value = 10.13588264717194
$ws.Range("E10").Value = 17.80151965644528
$ws.Range("F10").Value = 59.78810571157285
$ws.Range("G10").Value = 3.547611283709456
$ws.Range("I10").Value = 20.29074606183129
$ws.Range("M10").Value = 37.64702557419468
$ws.Range("B11").Value = 10.00807268459332
$ws.Range("D11").Value = 10.3786470747132
$ws.Range("E11").Value = 18.34925996050229
$ws.Range("F11").Value = 61.45801635186479
$ws.Range("G11").Value = 3.538341409436062
$ws.Range("I11").Value = 20.45317780383579
$ws.Range("M11").Value = 38.74754358735873
$ws.Range("B12").Value = 10.09501062114514
$ws.Range("D12").Value = 10.47065325504649
$ws.Range("E12").Value = 18.55250779540976
$ws.Range("F12").Value = 62.08637487379768
$ws.Range("G12").Value = 3.534872733183012
$ws.Range("I12").Value = 20.51540089589366
$ws.Range("M12").Value = 39.15803908702252
$ws.Range("B13").Value = 10.07635381368357
$ws.Range("D13").Value = 10.45083577663729
$ws.Range("E13").Value = 18.50892040269792
$ws.Range("F13").Value = 61.9512309020617
$ws.Range("G13").Value = 3.535617946664907
$ws.Range("I13").Value = 20.50196926222806
$ws.Range("M13").Value = 39.06991178903825
$ws.Range("B14").Value = 10.01525517817938
$ws.Range("D14").Value = 10.38621534376479
$ws.Range("E14").Value = 18.36606492907858
$ws.Range("F14").Value = 61.50979426892562
$ws.Range("G14").Value = 3.538055213481683
$ws.Range("I14").Value = 20.4582829783631
$ws.Range("M14").Value = 38.78144137721137
$ws.Range("B15").Value = 9.977635498377078
$ws.Range("D15").Value = 10.34664134993582
$ws.Range("E15").Value = 18.2780183861476
$ws.Range("F15").Value = 61.23886854237749
$ws.Range("G15").Value = 3.539553487360361
$ws.Range("I15").Value = 20.43161501578272
$ws.Range("M15").Value = 38.60392677523855
$ws.Range("B16").Value = 9.759095335508853
$ws.Range("D16").Value = 10.12003203502377
$ws.Range("E16").Value = 17.76514097356455
$ws.Range("F16").Value = 59.67843615606998
$ws.Range("G16").Value = 3.548222994943675
$ws.Range("I16").Value = 20.2802339407032
$ws.Range("M16").Value = 37.57423882388836
$ws.Range("B17").Value = 9.622426153955955
$ws.Range("D17").Value = 9.981228741174597
$ws.Range("E17").Value = 17.44309984048179
$ws.Range("F17").Value = 58.71449525269211
$ws.Range("G17").Value = 3.553617127467457
$ws.Range("I17").Value = 20.18870259652981
$ws.Range("M17").Value = 36.93161133884835
$ws.Range("B18").Value = 9.542876925554006
$ws.Range("D18").Value = 9.901496181491293
$ws.Range("E18").Value = 17.25515535281115
$ws.Range("F18").Value = 58.15776210779484
$ws.Range("G18").Value = 3.556747949167377
$ws.Range("I18").Value = 20.1365666165146
$ws.Range("M18").Value = 36.55802599161265
$ws.Range("B19").Value = 9.51578290480612
$ws.Range("D19").Value = 9.874520313005402
$ws.Range("E19").Value = 17.19105510021024
$ws.Range("F19").Value = 57.96888199614146
$ws.Range("G19").Value = 3.557812887204441
$ws.Range("I19").Value = 20.11900377472351
$ws.Range("M19").Value = 36.43086196293254
$ws.Range("B20").Value = 9.63707254700514
$ws.Range("D20").Value = 9.995994459455716
$ws.Range("E20").Value = 17.47766284106619
$ws.Range("F20").Value = 58.81734997619853
$ws.Range("G20").Value = 3.553039998155403
$ws.Range("I20").Value = 20.19839387449855
$ws.Range("M20").Value = 37.00043173027706
$ws.Range("B21").Value = 10.0332420401924
$ws.Range("D21").Value = 10.40519440938503
$ws.Range("E21").Value = 18.4081383201159
$ws.Range("F21").Value = 61.63956669046914
$ws.Range("G21").Value = 3.537338211072317
$ws.Range("I21").Value = 20.47109581198205
$ws.Range("M21").Value = 38.86634289399611
$ws.Range("B22").Value = 10.28347385383194
$ws.Range("D22").Value = 10.67304802433073
$ws.Range("E22").Value = 18.99195191723519
$ws.Range("F22").Value = 63.46057716023856
$ws.Range("G22").Value = 3.527318097349653
$ws.Range("I22").Value = 20.65345427360353
$ws.Range("M22").Value = 40.04935770363562
$ws.Range("B23").Value = 10.15072901362731
$ws.Range("D23").Value = 10.53007382104146
$ws.Range("E23").Value = 18.68258791871575
$ws.Range("F23").Value = 62.49094627580971
$ws.Range("G23").Value = 3.532644370319007
$ws.Range("I23").Value = 20.55576807109233
$ws.Range("M23").Value = 39.42134598372012
$ws.Range("B24").Value = 9.630453950012456
$ws.Range("D24").Value = 9.989318662757716
$ws.Range("E24").Value = 17.46204562428041
$ws.Range("F24").Value = 58.77085726516879
$ws.Range("G24").Value = 3.553300825817718
$ws.Range("I24").Value = 20.19401092746413
$ws.Range("M24").Value = 36.96933088197451
$ws.Range("B25").Value = 9.033546210045728
$ws.Range("D25").Value = 9.410361043430532
$ws.Range("E25").Value = 16.04188088321959
$ws.Range("F25").Value = 54.67278578349791
$ws.Range("G25").Value = 3.576667811113561
$ws.Range("I25").Value = 19.82355067188487
$ws.Range("M25").Value = 34.17405426644781
